$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of mortality-tracker data appended after the existing last row (151).
# Columns: A tracker_date, B report_date, C g_dead_total, D g_dead_child,
# E g_dead_women, F g_injured_total, G g_injured_child, H g_injured_women,
# I g_missing, J owb_dead_total, K owb_dead_children, L owb_injured_total, M source

$rows = @(
    @{ Row=152; A="07.03.2024"; B="06.03.2024"; C=30800; D=12300; E=8400;  F=72298; G=8663; H=6327; I=8000; J=424; K=113; L=4600; M="https://web.archive.org/web/20240307133101/https://www.aljazeera.com/news/longform/2023/10/9/israel-hamas-war-in-maps-and-charts-live-tracker" },
    @{ Row=153; A="08.03.2024"; B="08.03.2024"; C=30878; D=12300; E=8400;  F=72402; G=8663; H=6327; I=8000; J=424; K=113; L=4600; M="https://web.archive.org/web/20240308174707/https://www.aljazeera.com/news/longform/2023/10/9/israel-hamas-war-in-maps-and-charts-live-tracker" },
    @{ Row=154; A="09.03.2024"; B="08.03.2024"; C=30878; D=12300; E=8400;  F=72402; G=8663; H=6327; I=8000; J=424; K=113; L=4600; M="https://web.archive.org/web/20240309192855/https://www.aljazeera.com/news/longform/2023/10/9/israel-hamas-war-in-maps-and-charts-live-tracker" },
    @{ Row=155; A="10.03.2024"; B="10.03.2024"; C=31045; D=12300; E=8400;  F=72654; G=8663; H=6327; I=8000; J=425; K=113; L=4650; M="https://web.archive.org/web/20240310134417/https://www.aljazeera.com/news/longform/2023/10/9/israel-hamas-war-in-maps-and-charts-live-tracker" }
)

foreach ($r in $rows) {
    $row = $r.Row
    # Force A/B to be stored as plain text (they look like dates,
    # "DD.MM.YYYY", and would otherwise get auto-converted to date
    # serials by Excel's smart-entry parsing). Flip the format to Text
    # for the assignment, then restore General so no stray per-cell
    # number format sticks around on the saved cell.
    $ws.Cells.Item($row, 1).NumberFormat = "@"
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 1).NumberFormat = "General"
    $ws.Cells.Item($row, 2).NumberFormat = "General"
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
}

# Bring the frozen/split pane and selection up to date with the newly
# appended rows, mirroring how the sheet looked right after the last row
# (M155) was entered.
$win = $excel.ActiveWindow
$win.SplitRow = 128
$ws.Range("M155").Select() | Out-Null

